# Daily attendance processing - 2025-11-14 05:23:17
# Rotates the comma-separated "Recorded By" values in column G for the
# affected rows: the first name/email in the list is moved to the end.
# e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,6,10,11,12,13,14,15,17,18,19,20,21,22,24,26,28,29,32,36,37,38,39,40,41,43,44,45,46,47,48,50,52,54,55,58,62,63,64,65,66,67,69,70,71,72,73,74,76,78,83,84,85,86,87,90,92,93,94,96,99,101,109,110,111,112,113,116,118,119,120,122,125,127,135,136,137,138,139,142,144,145,146,148,151,153)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2
    $parts = $val -split ', '
    if ($parts.Length -gt 1) {
        $rotated = ($parts[1..($parts.Length - 1)] + $parts[0]) -join ', '
        $cell.Value2 = $rotated
    }
}
